$d = $word.ActiveDocument
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Consulta de producto por </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tipo_producto</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.codproductos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.nombre_producto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.tipo_producto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.fecha_ingreso</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.cantidad_producto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.precio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.ruc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.nombres,u.username</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">FROM </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd_systema.productos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> p </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">INNER JOIN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd_systema.proveedor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pro </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>on</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.idusuarios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.idusuarios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">INNER JOIN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd_systema.usuarios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> u </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>on</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.idusuarios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>u.idusuarios</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">ORDER BY </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.tipo_producto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Consulta de producto por </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rucproveedor</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.codproductos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.nombre_producto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.tipo_producto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.fecha_ingreso</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.cantidad_producto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.precio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.ruc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.nombres</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>u.username</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">FROM </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd_systema.productos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> p </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">INNER JOIN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd_systema.proveedor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pro </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>on</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.idusuarios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.idusuarios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">INNER JOIN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd_systema.usuarios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> u </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>on</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.idusuarios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>u.idusuarios</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">ORDER BY </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.ruc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>;</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.codproductos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.nombre_producto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.tipo_producto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.fecha_ingreso</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.cantidad_producto</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.precio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.ruc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.nombres</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.idProveedor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>u.username</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">FROM </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd_systema.productos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> p </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">INNER JOIN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd_systema.proveedor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pro </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>on</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.idProveedor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.idProveedor</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">INNER JOIN </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bd_systema.usuarios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> u </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>on</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p.idusuarios</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>u.idusuarios</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">ORDER BY </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pro.ruc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>;</w:t></w:r></w:p>'
$d.Content.InsertXML($xml)
